$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 28-32 (the rows containing the "-inf" placeholder values),
# which shifts all subsequent rows up by 5.
$ws.Range("A28:A32").EntireRow.Delete()
